# QCMI-BSA-offset.xlsx update:
# - fill in the computed offset values (col A,C,E,G,I,K,M) on row 2
# - give those "offset" cells a medium box border + centered, wrapped text
# - the header row goes back to the default (unstyled) cell style
# - row heights grow to fit the new formatting
# - column E is widened to fit its new (larger) numbers
# - the active selection moves to M2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- new offset values in row 2 ----
$ws.Range("A2").Value = 4961272.0199999996
$ws.Range("C2").Value = 14866726.869999999
$ws.Range("E2").Value = 24774372.23
$ws.Range("G2").Value = 34680641.270000003
$ws.Range("I2").Value = 44587621.439999998
$ws.Range("K2").Value = 54494404.460000001
$ws.Range("M2").Value = 64401754.57

# ---- header row (row 1) loses its (currently no-op) explicit style ----
$ws.Range("A1:N1").Style = "Normal"

# ---- build the new "boxed" style on A2, then stamp it onto the other ----
# ---- offset cells (C2,E2,G2,I2,K2,M2) without touching the fwhm cells ----
$boxed = $ws.Range("A2")
$boxed.Borders.Weight = -4138          # xlMedium
$boxed.VerticalAlignment = -4108       # xlCenter
$boxed.WrapText = $true
$boxed.Copy()
foreach ($col in @("C", "E", "G", "I", "K", "M")) {
    $ws.Range($col + "2").PasteSpecial(-4122)   # xlPasteFormats
}
$excel.CutCopyMode = 0

# ---- row heights grow slightly to accommodate the new formatting ----
$ws.Rows("1:2").RowHeight = 15.75

# ---- column E auto-sizes to fit the much larger numbers now in it ----
$ws.Columns("E").ColumnWidth = 10.7109375

# ---- move the selection to M2, matching the saved view state ----
[void]$ws.Range("M2").Select()

# ---- restore the workbook window position recorded on save ----
try {
    $win = $wb.Windows.Item(1)
    $win.Left = 13380
    $win.Top = 4005
} catch {
}
